$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    -501.8642238104799,
    -499.5404033536049,
    -497.21658289673,
    -494.8927624398549,
    -492.56894198298,
    -490.245121526105,
    -487.92130106923,
    -485.597480612355,
    -483.27366015548,
    -480.9498396986051,
    -478.62601924173,
    -476.3021987848551,
    -473.9783783279803,
    -471.6545578711052,
    -469.3307374142303,
    -467.0069169573552,
    -464.6830965004801,
    -462.3592760436053,
    -460.0354555867303,
    -457.7116351298553,
    -455.3878146729803,
    -453.0639942161054,
    -450.7401737592303,
    -448.4163533023553,
    -446.0925328454804,
    -443.7687123886055,
    -441.4448919317304,
    -439.1210714748555,
    -436.7972510179804,
    -434.4734305611054,
    -432.1496101042305,
    -429.8257896473555,
    -427.5019691904805,
    -425.1781487336055,
    -422.8543282767305,
    -420.5305078198554,
    -418.2066873629806,
    -415.8828669061056,
    -413.5590464492307,
    -411.2352259923557,
    -408.9114055354806,
    -406.5875850786057,
    -404.2637646217307,
    -401.9399441648557,
    -399.6161237079808,
    -397.2923032511057,
    -394.9684827942308,
    -392.6446623373558,
    -390.3208418804809,
    -387.9970214236058,
    -385.6732009667309,
    -383.3493805098558,
    -381.0255600529809,
    -378.701739596106,
    -376.3779191392309,
    -374.054098682356,
    -371.7302782254809,
    -369.4064577686059,
    -367.0826373117311,
    -364.758816854856,
    -362.4349963979811,
    -360.1111759411061,
    -357.7873554842312,
    -355.4635350273561,
    -353.1397145704811,
    -350.8158941136062,
    -348.4920736567311,
    -346.1682531998562,
    -343.8444327429811,
    -341.5206122861063,
    -339.1967918292312,
    -336.8729713723563,
    -334.5491509154813,
    -332.2253304586063,
    -329.9015100017313,
    -327.5776895448563,
    -325.2538690879814,
    -322.9300486311064,
    -320.6062281742314,
    -318.2824077173563,
    -315.9585872604815,
    -313.6347668036066,
    -311.3109463467315,
    -308.9871258898565,
    -306.6633054329815,
    -304.3394849761066,
    -302.0156645192316,
    -299.6918440623565,
    -297.3680236054817,
    -295.0442031486066,
    -292.7203826917316,
    -290.3965622348566,
    -288.0727417779816,
    -285.7489213211067,
    -283.4251008642317,
    -281.1012804073567,
    -278.7774599504817,
    -276.4536394936067,
    -274.1298190367318,
    -271.8059985798568
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
